$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '285.07'
$ws.Range("D2").ClearFormats()

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '2.13%'
$ws.Range("E2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '28.79'
$ws.Range("D3").ClearFormats()

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '5.12%'
$ws.Range("E3").ClearFormats()

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.954'
$ws.Range("D4").ClearFormats()

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '2.28%'
$ws.Range("E4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06675'
$ws.Range("D5").ClearFormats()

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '4.62%'
$ws.Range("E5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '7.326'
$ws.Range("D6").ClearFormats()

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '4.17%'
$ws.Range("E6").ClearFormats()

$ws.Range("B7").Value = 'GateToken'

$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.381'
$ws.Range("D7").ClearFormats()

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.37%'
$ws.Range("E7").ClearFormats()

$ws.Range("B8").Value = 'FTXToken'

$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.356'
$ws.Range("D8").ClearFormats()

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '2.25%'
$ws.Range("E8").ClearFormats()

$ws.Range("B9").Value = 'MXToken'

$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9373'
$ws.Range("D9").ClearFormats()

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '4.89%'
$ws.Range("E9").ClearFormats()

$ws.Range("B10").Value = 'WazirX'

$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1574'
$ws.Range("D10").ClearFormats()

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '3.09%'
$ws.Range("E10").ClearFormats()

$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'

$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06652'
$ws.Range("D11").ClearFormats()

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '16.96%'
$ws.Range("E11").ClearFormats()

$ws.Range("B12").Value = 'MandalaExchangeToken'

$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07677'
$ws.Range("D12").ClearFormats()

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '2.16%'
$ws.Range("E12").ClearFormats()

$ws.Range("B13").Value = 'BitrueCoin'

$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.02923'
$ws.Range("D13").ClearFormats()

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.03%'
$ws.Range("E13").ClearFormats()

$ws.Range("B14").Value = 'BitMartToken'

$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.08993'
$ws.Range("D14").ClearFormats()

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.07%'
$ws.Range("E14").ClearFormats()

$ws.Range("B15").Value = 'BitForexToken'

$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001597'
$ws.Range("D15").ClearFormats()

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.96%'
$ws.Range("E15").ClearFormats()

$ws.Range("B16").Value = 'CoinExToken'

$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.04478'
$ws.Range("D16").ClearFormats()

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.47%'
$ws.Range("E16").ClearFormats()

$ws.Range("B17").Value = 'One'

$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0006440'
$ws.Range("D17").ClearFormats()

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.50%'
$ws.Range("E17").ClearFormats()

$ws.Range("B18").Value = 'TigerCash'

$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006559'
$ws.Range("D18").ClearFormats()

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '7.61%'
$ws.Range("E18").ClearFormats()

$ws.Range("B19").Value = 'LEO'

$ws.Range("C19").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.491'
$ws.Range("D19").ClearFormats()

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.52%'
$ws.Range("E19").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.236'
$ws.Range("D20").ClearFormats()

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-2.04%'
$ws.Range("E20").ClearFormats()

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.3207'
$ws.Range("D21").ClearFormats()

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.91%'
$ws.Range("E21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1308'
$ws.Range("D22").ClearFormats()

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-3.17%'
$ws.Range("E22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.066'
$ws.Range("D23").ClearFormats()

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '4.25%'
$ws.Range("E23").ClearFormats()

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.22%'
$ws.Range("E24").ClearFormats()

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001178'
$ws.Range("D25").ClearFormats()

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.27%'
$ws.Range("E25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004483'
$ws.Range("D26").ClearFormats()

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '4.73%'
$ws.Range("E26").ClearFormats()

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001248'
$ws.Range("D27").ClearFormats()

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '5.87%'
$ws.Range("E27").ClearFormats()

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001613'
$ws.Range("D28").ClearFormats()

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '-2.40%'
$ws.Range("E28").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04186'
$ws.Range("D40").ClearFormats()

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '3.42%'
$ws.Range("E40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006746'
$ws.Range("D41").ClearFormats()

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '4.61%'
$ws.Range("E41").ClearFormats()

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-11.06%'
$ws.Range("E42").ClearFormats()

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002016'
$ws.Range("D43").ClearFormats()

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-7.82%'
$ws.Range("E43").ClearFormats()

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01219'
$ws.Range("D44").ClearFormats()

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '9.24%'
$ws.Range("E44").ClearFormats()

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005584'
$ws.Range("D45").ClearFormats()

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.45%'
$ws.Range("E45").ClearFormats()

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '20.74%'
$ws.Range("E46").ClearFormats()

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.01303'
$ws.Range("D47").ClearFormats()

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-29.51%'
$ws.Range("E47").ClearFormats()
